$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, week-covering dates) ---
$ws.Range("A8").Value = "Volume 30   Number  50"
$ws.Range("C9").Value = "Report Covering the Week  12/11/2023  Through  12/17/2023"

# --- Data table updates (rows 14-29) ---
# Cells whose style changes are handled with a Copy/PasteSpecial(xlPasteFormats=-4122)
# of formatting from a stable donor cell holding the exact target style, applied AFTER
# the value is written (paste-formats does not disturb the value already written).

$ws.Range("L14").Value = -25
$ws.Range("N14").Value = -76.923076923076
$ws.Range("C15").Value = 2
$ws.Range("C36").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$ws.Range("D15").Value = "'0"
$ws.Range("C34").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null
$ws.Range("E15").Value = "***.*"
$ws.Range("C34").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null
$ws.Range("F15").Value = 4
$ws.Range("H15").Value = 300
$ws.Range("I15").Value = 39
$ws.Range("K15").Value = 11.428571428571
$ws.Range("L15").Value = 11.428571428571
$ws.Range("M15").Value = 18.181818181818
$ws.Range("N15").Value = -58.064516129032
$ws.Range("C16").Value = 12
$ws.Range("E16").Value = -29.411764705882
$ws.Range("F16").Value = 53
$ws.Range("G16").Value = 50
$ws.Range("H16").Value = 6
$ws.Range("I16").Value = 539
$ws.Range("J16").Value = 683
$ws.Range("K16").Value = -21.083455344070
$ws.Range("L16").Value = 6.944444444444
$ws.Range("M16").Value = 8.016032064128
$ws.Range("N16").Value = -72.750252780586
$ws.Range("C17").Value = 14
$ws.Range("D17").Value = 14
$ws.Range("E17").Value = 0
$ws.Range("G17").Value = 69
$ws.Range("H17").Value = -4.347826086956
$ws.Range("I17").Value = 1039
$ws.Range("J17").Value = 1012
$ws.Range("K17").Value = 2.667984189723
$ws.Range("L17").Value = 24.431137724550
$ws.Range("M17").Value = 92.407407407407
$ws.Range("N17").Value = -21.938392186326
$ws.Range("C18").Value = 11
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 1000
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = 80.952380952380
$ws.Range("I18").Value = 445
$ws.Range("J18").Value = 405
$ws.Range("K18").Value = 9.876543209876
$ws.Range("L18").Value = 67.924528301886
$ws.Range("M18").Value = 68.560606060606
$ws.Range("N18").Value = -78.324403312226
$ws.Range("C19").Value = 14
$ws.Range("E19").Value = -12.5
$ws.Range("F19").Value = 43
$ws.Range("G19").Value = 51
$ws.Range("H19").Value = -15.686274509803
$ws.Range("I19").Value = 758
$ws.Range("J19").Value = 902
$ws.Range("K19").Value = -15.964523281596
$ws.Range("L19").Value = -3.684879288437
$ws.Range("M19").Value = 98.429319371727
$ws.Range("N19").Value = -2.319587628865
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = 12.5
$ws.Range("F20").Value = 27
$ws.Range("G20").Value = 33
$ws.Range("H20").Value = -18.181818181818
$ws.Range("I20").Value = 449
$ws.Range("J20").Value = 373
$ws.Range("K20").Value = 20.375335120643
$ws.Range("L20").Value = 172.121212121212
$ws.Range("M20").Value = 126.767676767677
$ws.Range("N20").Value = -63.673139158576
$ws.Range("C21").Value = 62
$ws.Range("D21").Value = 56
$ws.Range("E21").Value = 10.714285714285
$ws.Range("F21").Value = 231
$ws.Range("G21").Value = 225
$ws.Range("H21").Value = 2.666666666666
$ws.Range("I21").Value = 3284
$ws.Range("J21").Value = 3429
$ws.Range("K21").Value = -4.228638086905
$ws.Range("L21").Value = 25.775564917656
$ws.Range("M21").Value = 70.155440414507
$ws.Range("N21").Value = -56.399362719065
$ws.Range("C22").Value = 2
$ws.Range("C36").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null
$ws.Range("I22").Value = 39
$ws.Range("K22").Value = 5.405405405405
$ws.Range("L22").Value = 62.5
$ws.Range("M22").Value = -4.878048780487
$ws.Range("C23").Value = 2
$ws.Range("E23").Value = -33.333333333333
$ws.Range("F23").Value = 3
$ws.Range("H23").Value = -66.666666666666
$ws.Range("I23").Value = 74
$ws.Range("J23").Value = 81
$ws.Range("K23").Value = -8.641975308641
$ws.Range("L23").Value = -3.896103896103
$ws.Range("M23").Value = 54.166666666666
$ws.Range("C24").Value = 42
$ws.Range("D24").Value = 11
$ws.Range("E24").Value = 281.818181818182
$ws.Range("F24").Value = 135
$ws.Range("G24").Value = 118
$ws.Range("H24").Value = 14.406779661016
$ws.Range("I24").Value = 2002
$ws.Range("J24").Value = 1874
$ws.Range("K24").Value = 6.830309498399
$ws.Range("L24").Value = 48.736998514115
$ws.Range("M24").Value = 51.666666666666
$ws.Range("C25").Value = 22
$ws.Range("D25").Value = 17
$ws.Range("E25").Value = 29.411764705882
$ws.Range("F25").Value = 86
$ws.Range("G25").Value = 65
$ws.Range("H25").Value = 32.307692307692
$ws.Range("I25").Value = 1189
$ws.Range("J25").Value = 1203
$ws.Range("K25").Value = -1.163757273482
$ws.Range("L25").Value = 16.11328125
$ws.Range("M25").Value = -0.418760469011
$ws.Range("C26").Value = 3
$ws.Range("C36").Copy() | Out-Null
$ws.Range("C26").PasteSpecial(-4122) | Out-Null
$ws.Range("D26").Value = "'0"
$ws.Range("C34").Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4122) | Out-Null
$ws.Range("E26").Value = "***.*"
$ws.Range("C34").Copy() | Out-Null
$ws.Range("E26").PasteSpecial(-4122) | Out-Null
$ws.Range("F26").Value = 5
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 400
$ws.Range("I26").Value = 66
$ws.Range("K26").Value = 17.857142857142
$ws.Range("L26").Value = 15.789473684210
$ws.Range("C27").Value = 3
$ws.Range("E27").Value = 50
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 16.666666666666
$ws.Range("I27").Value = 128
$ws.Range("J27").Value = 111
$ws.Range("K27").Value = 15.315315315315
$ws.Range("L27").Value = 29.292929292929
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = -60
$ws.Range("L28").Value = -37.209302325581
$ws.Range("N28").Value = -68.965517241379
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 5
$ws.Range("H29").Value = -60
$ws.Range("L29").Value = -36.231884057971
$ws.Range("N29").Value = -72.151898734177
